# Applies the "Add validation.py + various corrections" edit to the
# BIGCC LCI workbook: refreshes the electrical-efficiency assumption used
# by the two "electricity production, at biomass-fired IGCC power plant"
# datasets (no-CCS block rows 1-28, with-CCS block rows 30-56), fixes a
# typo in the no-CCS comment, extends both process comments with notes
# about the efficiency assumption, and re-enters the source citation on
# the with-CCS block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Block 1 (rows 1-28): "electricity production, at biomass-fired IGCC
# power plant" (no CCS)
# ---------------------------------------------------------------------

# B9 ("source") had an accidental alignment style applied to it -
# clear it back to the default Normal style.
$ws.Range("B9").Style = "Normal"

# B10 ("comment") - extend with the note about the 19% -> 35% efficiency
# scale-up, and apply the (previously-stray) alignment style that used
# to live on B9.
$b10Text = @"
The study is originally tailored to portuguese conditions, using Portugal-grown eucalyptus as biommas feedstock. The plant has a power output of 12.5 MW.
The construciton, maintenance and operation of the plant are not from this study, but from Volkart et la, 2013. The required amounts for those have been adjusted to 12.5 MW.
All the inventory data regarding the RFB direct gasification process were obtained from experiments carried out in a pilot-scale gasification installation running at the University of Aveiro, Portugal (Pio et al., 2017). However, data are scaled up in order to obtain the same power output of the EG-CRC power plant. At the pilot scale, the RFB is introduced in the gasification chamber by means of a screw feeder and is converted at an average temperature of 785 °C in a bubbling fluidized bed reactor of 80 kWth, operated at atmospheric pressure and under auto-thermal regime, thus, direct gasification using atmospheric air. Note that this study considers an electricial efficiency of 19%, which is very low compared to the rest of the literature. We scaled up this efficiency to 35%, a central estimate when considering other studies: Carpentieri et al. (2005), Puy et al. (2010), Siegl et al. (2011), Guest et al. (2011), Steubing et al. (2011), Nguyen et al. (2013), Jäppinen et al. (2014), Wang et al. (2014), Paengjuntuek et al. (2015), Cambero et al. (2015), Yang et al. (2018) and Zang et al. (2020).
"@
$ws.Range("B10").Value = $b10Text
$ws.Range("B10").HorizontalAlignment = 1
$ws.Range("B10").VerticalAlignment = -4107

# B16 (wood input, kg/kWh) - was a hard-coded constant (1.027); now
# derived from the 35% electrical efficiency assumption: (3.6 / 35%) / 20
$ws.Range("B16").Formula = "=(3.6/35%)/20"
$ws.Range("B16").NumberFormat = "0.000"

# ---------------------------------------------------------------------
# Block 2 (rows 30-56): "electricity production, at biomass-fired IGCC
# power plant, pre, pipeline 200km, storage 1000m" (with CCS)
# ---------------------------------------------------------------------

# B36 ("source") - re-enter the citation (functionally identical text).
$citationText = "Andrei Briones-Hidrovo, José Copa, Luís A.C. Tarelho, Cátia Gonçalves, Tamíris Pacheco da Costa, Ana Cláudia Dias, Environmental and energy performance of residual forest biomass for electricity generation: Gasification vs. combustion, Journal of Cleaner Production, Volume 289, 2021, 125680, ISSN 0959-6526, https://doi.org/10.1016/j.jclepro.2020.125680."
$ws.Range("B36").Value = $citationText

# B37 ("comment") - fix "construciton" -> "construction", extend with the
# same efficiency note, and add the CCS efficiency-penalty note. Also
# apply the alignment style (matches B10's).
$b37Text = @"
The study is originally tailored to portuguese conditions, using Portugal-grown eucalyptus as biommas feedstock. The plant has a power output of 12.5 MW.
Also, the original study does not include CCS. This is added from Volkart et al, 2013.
The construction, maintenance and operation of the plant are not from this study, but from Volkart et la, 2013. The required amounts for those have been adjusted to 12.5 MW.
All the inventory data regarding the RFB direct gasification process were obtained from experiments carried out in a pilot-scale gasification installation running at the University of Aveiro, Portugal (Pio et al., 2017). However, data are scaled up in order to obtain the same power output of the EG-CRC power plant. At the pilot scale, the RFB is introduced in the gasification chamber by means of a screw feeder and is converted at an average temperature of 785 °C in a bubbling fluidized bed reactor of 80 kWth, operated at atmospheric pressure and under auto-thermal regime, thus, direct gasification using atmospheric air. Note that this study considers an electricial efficiency of 19%, which is very low compared to the rest of the literature. We scaled up this efficiency to 35%, a central estimate when considering other studies: Carpentieri et al. (2005), Puy et al. (2010), Siegl et al. (2011), Guest et al. (2011), Steubing et al. (2011), Nguyen et al. (2013), Jäppinen et al. (2014), Wang et al. (2014), Paengjuntuek et al. (2015), Cambero et al. (2015), Yang et al. (2018) and Zang et al. (2020). And then, we reduce the efficiency by approximately 10% (ie.., 30% eff.) to reflect the addendum of a CCS unit.
"@
$ws.Range("B37").Value = $b37Text
$ws.Range("B37").HorizontalAlignment = 1
$ws.Range("B37").VerticalAlignment = -4107

# B43 (wood input, kg/kWh) - same as B16 but for the 30% (CCS-penalised)
# electrical efficiency assumption: (3.6 / 30%) / 20
$ws.Range("B43").Formula = "=(3.6/30%)/20"
$ws.Range("B43").NumberFormat = "0.000"

[void]$wb.Save()
